$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 399, shifting the existing data (old rows
# 399-496) down to rows 400-497, and fill the new row 399 with the
# latest week's price record for Rabanito at Vega Central Mapocho de
# Santiago.
$ws.Rows(399).Insert()

$ws.Range("A399").Value = 9
$ws.Range("B399").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C399").Value = "Metropolitana"
$ws.Range("D399").Value = 45204
$ws.Range("E399").Value = 13
$ws.Range("F399").Value = 300000001
$ws.Range("G399").Value = "Rabanito"
$ws.Range("H399").Value = "Sin especificar"
$ws.Range("I399").Value = "Primera"
$ws.Range("J399").Value = 7000
$ws.Range("K399").Value = 3000
$ws.Range("L399").Value = 3000
$ws.Range("M399").Value = 3000
$ws.Range("N399").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O399").Value = "Provincia de Chacabuco"
$ws.Range("P399").Value = 30
$ws.Range("Q399").Value = 100
$ws.Range("R399").Value = "Hortaliza"
